$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, pushing existing rows 6-35 down to 7-36.
$ws.Rows.Item(6).Insert()

# Populate the new row with the "start" entries (matches style inherited from row above).
$ws.Range("A6").Value = "start"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "start"

# Update the selection to match the target state.
[void]$ws.Range("H11").Select()
